# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect the new report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-03 19:19:18"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-03 19:19:14"
$zhcn.Range("K2").Value = "2016-09-03 19:19:32"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-03 19:19:18"
$dede.Range("K2").Value = "2016-09-03 19:19:39"
